$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New classes appended to the schedule (rows 25-37). The worksheet already
# carries pre-formatted (but empty) placeholder cells for column A (date,
# style s=4) and column D (hyperlink, style s=5) down through row 39 - except
# for rows 28, 33 and 37 which don't exist yet at all, so their date cell's
# number format is copied from an existing date cell (A10) before the value
# is written. Columns B/C fall back to the column-level default style
# (s=2 / s=3) when a cell has no explicit style of its own, matching the
# existing rows above.
# ---------------------------------------------------------------------------

# --- 4/9/2020 - Framework para escalar Scrum --------------------------------
$ws.Range("A25").Value = 44078
$ws.Range("B25").Value = "Framework para escalar Scrum"
$ws.Range("C25").Value = 1
$ws.Hyperlinks.Add($ws.Range("D25"), "https://youtu.be/tx-tgq30vco")

$ws.Range("C26").Value = 2
$ws.Hyperlinks.Add($ws.Range("D26"), "https://youtu.be/bqUNI8AFqKA")

# --- 18/9/2020 - Testing -----------------------------------------------------
$ws.Range("A10").Copy($ws.Range("A28"))
$ws.Range("A28").Value = 44092
$ws.Range("B28").Value = "Testing"
$ws.Range("C28").Value = 1
$ws.Hyperlinks.Add($ws.Range("D28"), "https://youtu.be/-6vAPyi28OU")

$ws.Range("C29").Value = 2
$ws.Hyperlinks.Add($ws.Range("D29"), "https://youtu.be/rgRES2s5a_8")

# --- 25/9/2020 - Testing Agile en contexto ----------------------------------
$ws.Range("A10").Copy($ws.Range("A31"))
$ws.Range("A31").Value = 44099
$ws.Range("B31").Value = "Testing Agile en contexto"
$ws.Range("C31").Value = 1
$ws.Hyperlinks.Add($ws.Range("D31"), "https://youtu.be/hKoJBlhxuN8")

# --- 2/10/2020 - PPQA ---------------------------------------------------------
$ws.Range("A10").Copy($ws.Range("A33"))
$ws.Range("A33").Value = 44106
$ws.Range("B33").Value = "PPQA"
$ws.Range("C33").Value = 1
$ws.Hyperlinks.Add($ws.Range("D33"), "https://youtu.be/6GVTACucmRU")

# --- 9/10/2020 - Lean y Kanban ------------------------------------------------
$ws.Range("A10").Copy($ws.Range("A35"))
$ws.Range("A35").Value = 44113
$ws.Range("B35").Value = "Lean y Kanban"
$ws.Range("C35").Value = 1
$ws.Hyperlinks.Add($ws.Range("D35"), "https://youtu.be/opDHtLul3Hk")

# --- 23/10/2020 - Métricas -----------------------------------------------------
$ws.Range("A10").Copy($ws.Range("A37"))
$ws.Range("A37").Value = 44127
$ws.Range("B37").Value = "Métricas"
$ws.Range("C37").Value = 1
$ws.Hyperlinks.Add($ws.Range("D37"), "https://youtu.be/zoKbW8MlbsQ")
